# Processed remaining UWPR 2017 ETNP samples:
# add a new "PEAKS DB % mod peptides" column (inserted before the existing
# "% PSM Carbamidomethylated" column) and fill in the computed percentages
# for the newly-processed MCLANE/trap samples.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column at L; everything from the old L ("% PSM
# Carbamidomethylated") rightwards shifts over to M.. and inherits the
# per-row formatting automatically.
$ws.Columns("L").Insert()

# Header for the freshly inserted column.
$ws.Range("L1").Value2 = "PEAKS DB % mod peptides"

# Newly computed "% mod peptides" values for the processed samples.
$ws.Range("L29").Value2 = 0.53968253968253899
$ws.Range("L30").Value2 = 0.47435669920141899
$ws.Range("L31").Value2 = 0.42578125
$ws.Range("L32").Value2 = 0.29838709677419301
$ws.Range("L33").Value2 = 0.734513274336283
$ws.Range("L34").Value2 = 0.43161814488828698
$ws.Range("L35").Value2 = 0.64808481532147699

# Widen column B a touch so the longer sample names still fit.
$ws.Columns("B").ColumnWidth = 24.5

# Leave the selection where the author's last edit was.
$ws.Range("L34").Select()
